$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the sentence that was split by a stray "_GoBack" bookmark into a
#    single run (the bookmark disappears from this location).
# ---------------------------------------------------------------------------
$sentence = "written to provide a general SPM analysis pipeline of fMRI data, individual functions and the overall structure might be useful in developing new pipelines and possibly a more universally applicable script."
$d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a reviewer comment anchored on the "References" heading.
# ---------------------------------------------------------------------------
$word.UserName = "Anna V"
$word.UserInitials = "AV"

$headingRange = $d.Content
$headingRange.Find.Execute("References", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Comments.Add($headingRange, "Add wfu pickatlas aal reference") | Out-Null

# ---------------------------------------------------------------------------
# 3) Split the Ashburner reference run and drop a new "_GoBack" bookmark at
#    that edit point (mirrors Word automatically relocating the bookmark to
#    the most recent edit location).
# ---------------------------------------------------------------------------
$splitAnchor = $d.Content
$splitAnchor.Find.Execute("Kilner, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $splitAnchor.End
$goBackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
